$d = $word.ActiveDocument

# Update the date/day heading paragraph.
$d.Content.Find.Execute("2025-02-23 Sunday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-02-24 Monday", 2)

# Update the division-problem table cells by position, since several
# new values collide with other old values elsewhere in the table
# (a global text Find/Replace would mis-fire on those).
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "83÷3=" },
    @{ Row = 1;  Col = 2; Text = "88÷7=" },
    @{ Row = 1;  Col = 3; Text = "46÷9=" },
    @{ Row = 1;  Col = 4; Text = "17÷7=" },
    @{ Row = 1;  Col = 5; Text = "67÷7=" },

    @{ Row = 5;  Col = 1; Text = "57÷9=" },
    @{ Row = 5;  Col = 2; Text = "75÷4=" },
    @{ Row = 5;  Col = 3; Text = "88÷8=" },
    @{ Row = 5;  Col = 4; Text = "20÷2=" },
    @{ Row = 5;  Col = 5; Text = "70÷9=" },

    @{ Row = 9;  Col = 1; Text = "69÷3=" },
    @{ Row = 9;  Col = 2; Text = "18÷3=" },
    @{ Row = 9;  Col = 3; Text = "94÷4=" },
    @{ Row = 9;  Col = 4; Text = "81÷9=" },
    @{ Row = 9;  Col = 5; Text = "67÷7=" },

    @{ Row = 13; Col = 1; Text = "41÷4=" },
    @{ Row = 13; Col = 2; Text = "95÷7=" },
    @{ Row = 13; Col = 3; Text = "19÷4=" },
    @{ Row = 13; Col = 4; Text = "39÷6=" },
    @{ Row = 13; Col = 5; Text = "98÷8=" },

    @{ Row = 17; Col = 1; Text = "83÷7=" },
    @{ Row = 17; Col = 2; Text = "15÷9=" },
    @{ Row = 17; Col = 3; Text = "67÷8=" },
    @{ Row = 17; Col = 4; Text = "57÷3=" },
    @{ Row = 17; Col = 5; Text = "91÷5=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cellRange = $cell.Range
    [void]$cellRange.MoveEnd(1, -1)  # exclude trailing end-of-cell marker
    $cellRange.Text = $u.Text
}
